# Update the answer table in place: replace the division-equation text
# in specific cells of the (single) table with the new values. Using
# Table.Cell(row, col).Range.Text preserves the existing run/paragraph
# formatting (font, size, justification) while swapping only the text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "75÷9=8, 3" },
    @{ Row = 1;  Col = 2; Text = "50÷5=10, 0" },
    @{ Row = 1;  Col = 3; Text = "14÷3=4, 2" },
    @{ Row = 1;  Col = 4; Text = "29÷2=14, 1" },
    @{ Row = 1;  Col = 5; Text = "59÷7=8, 3" },

    @{ Row = 5;  Col = 1; Text = "32÷6=5, 2" },
    @{ Row = 5;  Col = 2; Text = "33÷8=4, 1" },
    @{ Row = 5;  Col = 3; Text = "34÷5=6, 4" },
    @{ Row = 5;  Col = 4; Text = "80÷5=16, 0" },
    @{ Row = 5;  Col = 5; Text = "25÷5=5, 0" },

    @{ Row = 9;  Col = 1; Text = "53÷5=10, 3" },
    @{ Row = 9;  Col = 2; Text = "46÷5=9, 1" },
    @{ Row = 9;  Col = 3; Text = "89÷7=12, 5" },
    @{ Row = 9;  Col = 4; Text = "75÷9=8, 3" },
    @{ Row = 9;  Col = 5; Text = "83÷3=27, 2" },

    @{ Row = 13; Col = 1; Text = "58÷5=11, 3" },
    @{ Row = 13; Col = 2; Text = "51÷5=10, 1" },
    @{ Row = 13; Col = 3; Text = "94÷2=47, 0" },
    @{ Row = 13; Col = 4; Text = "95÷8=11, 7" },
    @{ Row = 13; Col = 5; Text = "54÷2=27, 0" },

    @{ Row = 17; Col = 1; Text = "94÷4=23, 2" },
    @{ Row = 17; Col = 2; Text = "60÷6=10, 0" },
    @{ Row = 17; Col = 3; Text = "88÷8=11, 0" },
    @{ Row = 17; Col = 4; Text = "12÷5=2, 2" },
    @{ Row = 17; Col = 5; Text = "43÷8=5, 3" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}

Write-Output "Updated $($updates.Count) cells"
